# New crime data collected — roll the weekly CompStat report forward one
# week (Volume 30, Number 42 -> 43; week of 10/16-10/22/2023 -> 10/23-10/29/2023)
# and refresh the underlying precinct crime-complaint figures (rows 15-30).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Report header / volume number (A8) and week-covered banner (C9) ---
# These source cells hold multi-run rich text in the original file; COM
# collapses that to a single run on write, which is the best a COM-level
# edit can do, but the visible text ends up identical to the authored diff.
$ws.Range("A8").Value = "Volume 30   Number  43"
$ws.Range("C9").Value = "Report Covering the Week  10/23/2023  Through  10/29/2023"

# --- Row 15 (Rape) ---
$ws.Range("F15").Value = 2
# G15/H15 flip from numeric to the "no data" text markers ("0" / "***.*");
# copy from row 14's already-styled source cells so we reuse the existing
# cell style (14) instead of minting a new one via NumberFormat changes.
$ws.Range("F14").Copy($ws.Range("G15"))
$ws.Range("E14").Copy($ws.Range("H15"))
$ws.Range("N15").Value = -67.647058823529

# --- Row 16 (Robbery) ---
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 25
$ws.Range("F16").Value = 13
$ws.Range("G16").Value = 12
$ws.Range("H16").Value = 8.333333333333
$ws.Range("I16").Value = 154
$ws.Range("J16").Value = 147
$ws.Range("K16").Value = 4.761904761904
$ws.Range("L16").Value = 43.925233644859
$ws.Range("M16").Value = -50.482315112540
$ws.Range("N16").Value = -85.701021355617

# --- Row 17 (Fel. Assault) ---
$ws.Range("C17").Value = 7
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = 16.666666666666
$ws.Range("F17").Value = 25
$ws.Range("G17").Value = 16
$ws.Range("H17").Value = 56.25
$ws.Range("I17").Value = 254
$ws.Range("J17").Value = 235
$ws.Range("K17").Value = 8.085106382978
$ws.Range("L17").Value = 39.560439560439
$ws.Range("M17").Value = 51.190476190476
$ws.Range("N17").Value = -49.098196392785

# --- Row 18 (Burglary) ---
$ws.Range("C18").Value = 6
$ws.Range("D18").Value = 8
$ws.Range("E18").Value = -25
$ws.Range("F18").Value = 24
$ws.Range("G18").Value = 30
$ws.Range("H18").Value = -20
$ws.Range("I18").Value = 232
$ws.Range("J18").Value = 252
$ws.Range("K18").Value = -7.936507936507
$ws.Range("L18").Value = 25.405405405405
$ws.Range("M18").Value = -40.512820512820
$ws.Range("N18").Value = -79.450841452612

# --- Row 19 (Gr. Larceny) ---
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 16
$ws.Range("E19").Value = -37.5
$ws.Range("F19").Value = 49
$ws.Range("G19").Value = 57
$ws.Range("H19").Value = -14.035087719298
$ws.Range("I19").Value = 627
$ws.Range("J19").Value = 558
$ws.Range("K19").Value = 12.365591397849
$ws.Range("L19").Value = 54.433497536945
$ws.Range("M19").Value = 52.926829268292
$ws.Range("N19").Value = 39.024390243902

# --- Row 20 (G.L.A.) ---
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = -33.333333333333
$ws.Range("F20").Value = 11
$ws.Range("G20").Value = 14
$ws.Range("H20").Value = -21.428571428571
$ws.Range("I20").Value = 149
$ws.Range("J20").Value = 150
$ws.Range("K20").Value = -0.666666666666
$ws.Range("L20").Value = 14.615384615384
$ws.Range("M20").Value = 10.370370370370
$ws.Range("N20").Value = -80.471821756225

# --- Row 21 (TOTAL) ---
$ws.Range("C21").Value = 30
$ws.Range("D21").Value = 37
$ws.Range("E21").Value = -18.918918918918
$ws.Range("F21").Value = 124
$ws.Range("G21").Value = 129
$ws.Range("H21").Value = -3.875968992248
$ws.Range("I21").Value = 1432
$ws.Range("J21").Value = 1357
$ws.Range("K21").Value = 5.526897568165
$ws.Range("L21").Value = 38.894277400582
$ws.Range("M21").Value = 0.632466619817
$ws.Range("N21").Value = -63.947633434038

# --- Row 22 (Transit) ---
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -48.648648648648

# --- Row 23 (Housing) ---
$ws.Range("C23").Value = 5
$ws.Range("D23").Value = 5
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 18
$ws.Range("G23").Value = 13
$ws.Range("H23").Value = 38.461538461538
$ws.Range("I23").Value = 166
$ws.Range("J23").Value = 131
$ws.Range("K23").Value = 26.717557251908
$ws.Range("L23").Value = 20.289855072463
$ws.Range("M23").Value = 40.677966101694

# --- Row 24 (Petit Larceny) ---
$ws.Range("C24").Value = 17
$ws.Range("D24").Value = 30
$ws.Range("E24").Value = -43.333333333333
$ws.Range("F24").Value = 78
$ws.Range("G24").Value = 86
$ws.Range("H24").Value = -9.302325581395
$ws.Range("I24").Value = 880
$ws.Range("J24").Value = 1017
$ws.Range("K24").Value = -13.470993117010
$ws.Range("L24").Value = 9.862671660424
$ws.Range("M24").Value = -14.396887159533

# --- Row 25 (Misd. Assault) ---
$ws.Range("C25").Value = 16
$ws.Range("D25").Value = 9
$ws.Range("E25").Value = 77.777777777777
$ws.Range("F25").Value = 34
$ws.Range("G25").Value = 35
$ws.Range("H25").Value = -2.857142857142
$ws.Range("I25").Value = 413
$ws.Range("J25").Value = 409
$ws.Range("K25").Value = 0.977995110024
$ws.Range("L25").Value = 33.656957928802
$ws.Range("M25").Value = 2.227722772277

# --- Row 26 (UCR Rape*) ---
$ws.Range("F26").Value = 2
# G26/H26 flip to the "no data" text markers, same trick as row 15.
$ws.Range("F14").Copy($ws.Range("G26"))
$ws.Range("E14").Copy($ws.Range("H26"))

# --- Row 27 (Other Sex Crimes) ---
$ws.Range("D27").Value = 3
$ws.Range("E27").Value = -33.333333333333
$ws.Range("F27").Value = 4
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 48
$ws.Range("J27").Value = 40
$ws.Range("K27").Value = 20
$ws.Range("L27").Value = -22.580645161290

# --- Row 28 (Shooting Vic.) ---
$ws.Range("M28").Value = -14.285714285714
$ws.Range("N28").Value = -84.210526315789

# --- Row 29 (Shooting Inc.) ---
$ws.Range("M29").Value = -30.769230769230
$ws.Range("N29").Value = -87.142857142857

# --- Row 30 (Hate Crimes) ---
$ws.Range("D30").Value = 2
# F30 flips from numeric to the "0" text marker.
$ws.Range("F14").Copy($ws.Range("F30"))
$ws.Range("G30").Value = 3
$ws.Range("H30").Value = -100
$ws.Range("J30").Value = 28
$ws.Range("K30").Value = -50
